# Auto-generated edit script: updates currentAveragePrice/profit columns (H:N)
# across multiple class sheets per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 99.5
$ws.Cells.Item(5, 9).Value = 107.75
$ws.Cells.Item(5, 10).Value = 74.75
$ws.Cells.Item(5, 11).Value = 107.75
$ws.Cells.Item(5, 12).Value = 74.75
$ws.Cells.Item(5, 13).Value = 7.25
$ws.Cells.Item(5, 14).Value = -304.75
$ws.Cells.Item(17, 8).Value = 481259.97
$ws.Cells.Item(17, 10).Value = 481259.97
$ws.Cells.Item(17, 12).Value = 1443779.91
$ws.Cells.Item(17, 14).Value = -1444115.91
$ws.Cells.Item(94, 8).Value = 799
$ws.Cells.Item(94, 9).Value = 799
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 799
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -348
$ws.Cells.Item(94, 14).ClearContents() | Out-Null
$ws.Cells.Item(106, 8).Value = 3775.5557
$ws.Cells.Item(106, 9).Value = 3720
$ws.Cells.Item(106, 11).Value = 3720
$ws.Cells.Item(106, 13).Value = -3089
$ws.Cells.Item(132, 8).Value = 2085.2778
$ws.Cells.Item(132, 9).Value = 2329
$ws.Cells.Item(132, 10).Value = 866.6667
$ws.Cells.Item(132, 11).Value = 6987
$ws.Cells.Item(132, 12).Value = 2600.0001
$ws.Cells.Item(132, 13).Value = -4457
$ws.Cells.Item(132, 14).Value = -7660.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1585.3334
$ws.Cells.Item(45, 9).Value = 830.5
$ws.Cells.Item(45, 10).Value = 3095
$ws.Cells.Item(45, 11).Value = 830.5
$ws.Cells.Item(45, 12).Value = 3095
$ws.Cells.Item(45, 13).Value = -453.5
$ws.Cells.Item(45, 14).Value = -3849
$ws.Cells.Item(110, 8).Value = 706.3570999999999
$ws.Cells.Item(110, 9).Value = 682.7692
$ws.Cells.Item(110, 10).Value = 1013
$ws.Cells.Item(110, 11).Value = 682.7692
$ws.Cells.Item(110, 12).Value = 1013
$ws.Cells.Item(110, 13).Value = 1362.2308
$ws.Cells.Item(110, 14).Value = -5103

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 14).ClearContents() | Out-Null
$ws.Cells.Item(105, 8).Value = 9980.869000000001
$ws.Cells.Item(105, 9).Value = 8920
$ws.Cells.Item(105, 10).Value = 13800
$ws.Cells.Item(105, 11).Value = 8920
$ws.Cells.Item(105, 12).Value = 13800
$ws.Cells.Item(105, 13).Value = -7173
$ws.Cells.Item(105, 14).Value = -17294
$ws.Cells.Item(134, 8).Value = 24841.592
$ws.Cells.Item(134, 9).Value = 2053.8386
$ws.Cells.Item(134, 11).Value = 6161.5158
$ws.Cells.Item(134, 13).Value = -3626.5158

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 2678.125
$ws.Cells.Item(86, 9).Value = 2618.75
$ws.Cells.Item(86, 10).Value = 2737.5
$ws.Cells.Item(86, 11).Value = 2618.75
$ws.Cells.Item(86, 12).Value = 2737.5
$ws.Cells.Item(86, 13).Value = -1495.75
$ws.Cells.Item(86, 14).Value = -4983.5
$ws.Cells.Item(89, 8).Value = 2678.125
$ws.Cells.Item(89, 9).Value = 2618.75
$ws.Cells.Item(89, 10).Value = 2737.5
$ws.Cells.Item(89, 11).Value = 13093.75
$ws.Cells.Item(89, 12).Value = 13687.5
$ws.Cells.Item(89, 13).Value = -7477.75
$ws.Cells.Item(89, 14).Value = -24919.5
$ws.Cells.Item(103, 8).Value = 29250
$ws.Cells.Item(103, 9).Value = 3000
$ws.Cells.Item(103, 10).Value = 38000
$ws.Cells.Item(103, 11).Value = 3000
$ws.Cells.Item(103, 12).Value = 38000
$ws.Cells.Item(103, 13).Value = -1828
$ws.Cells.Item(103, 14).Value = -40344
$ws.Cells.Item(132, 8).Value = 2132.6843
$ws.Cells.Item(132, 9).Value = 1301.7693
$ws.Cells.Item(132, 10).Value = 3933
$ws.Cells.Item(132, 11).Value = 3905.3079
$ws.Cells.Item(132, 12).Value = 11799
$ws.Cells.Item(132, 13).Value = -1375.3079
$ws.Cells.Item(132, 14).Value = -16859

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4451.6875
$ws.Cells.Item(70, 9).Value = 4372.5713
$ws.Cells.Item(70, 10).Value = 4513.222
$ws.Cells.Item(70, 11).Value = 4372.5713
$ws.Cells.Item(70, 12).Value = 4513.222
$ws.Cells.Item(70, 13).Value = -4102.5713
$ws.Cells.Item(70, 14).Value = -5053.222
$ws.Cells.Item(73, 8).Value = 4451.6875
$ws.Cells.Item(73, 9).Value = 4372.5713
$ws.Cells.Item(73, 10).Value = 4513.222
$ws.Cells.Item(73, 11).Value = 4372.5713
$ws.Cells.Item(73, 12).Value = 4513.222
$ws.Cells.Item(73, 13).Value = -3436.5713
$ws.Cells.Item(73, 14).Value = -6385.222
$ws.Cells.Item(80, 8).Value = 3583.3333
$ws.Cells.Item(80, 9).Value = 3500
$ws.Cells.Item(80, 10).Value = 4000
$ws.Cells.Item(80, 11).Value = 3500
$ws.Cells.Item(80, 12).Value = 4000
$ws.Cells.Item(80, 13).Value = -2502
$ws.Cells.Item(80, 14).Value = -5996
$ws.Cells.Item(83, 8).Value = 3583.3333
$ws.Cells.Item(83, 9).Value = 3500
$ws.Cells.Item(83, 10).Value = 4000
$ws.Cells.Item(83, 11).Value = 17500
$ws.Cells.Item(83, 12).Value = 20000
$ws.Cells.Item(83, 13).Value = -12508
$ws.Cells.Item(83, 14).Value = -29984
$ws.Cells.Item(102, 8).Value = 2451.6316
$ws.Cells.Item(102, 9).Value = 1760.0769
$ws.Cells.Item(102, 10).Value = 3950
$ws.Cells.Item(102, 11).Value = 1760.0769
$ws.Cells.Item(102, 12).Value = 3950
$ws.Cells.Item(102, 13).Value = -138.0769
$ws.Cells.Item(102, 14).Value = -7194
$ws.Cells.Item(132, 8).Value = 2690.5
$ws.Cells.Item(132, 9).Value = 2139.8096
$ws.Cells.Item(132, 11).Value = 6419.4288
$ws.Cells.Item(132, 13).Value = -3889.4288

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 18521544
$ws.Cells.Item(7, 10).Value = 55556810
$ws.Cells.Item(7, 12).Value = 55556810
$ws.Cells.Item(7, 14).Value = -55557034
$ws.Cells.Item(40, 8).Value = 1011561
$ws.Cells.Item(40, 9).Value = 2021512
$ws.Cells.Item(40, 11).Value = 2021512
$ws.Cells.Item(40, 13).Value = -2021376
$ws.Cells.Item(100, 8).Value = 2000
$ws.Cells.Item(100, 9).Value = 1975.125
$ws.Cells.Item(100, 10).Value = 2066.3333
$ws.Cells.Item(100, 11).Value = 1975.125
$ws.Cells.Item(100, 12).Value = 2066.3333
$ws.Cells.Item(100, 13).Value = -1434.125
$ws.Cells.Item(100, 14).Value = -3148.3333
$ws.Cells.Item(126, 8).Value = 18521544
$ws.Cells.Item(126, 10).Value = 55556810
$ws.Cells.Item(126, 12).Value = 166670430
$ws.Cells.Item(126, 14).Value = -166675370
$ws.Cells.Item(130, 8).Value = 0
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 14).ClearContents() | Out-Null
$ws.Cells.Item(132, 8).Value = 1717.1482
$ws.Cells.Item(132, 9).Value = 1511
$ws.Cells.Item(132, 10).Value = 3366.3333
$ws.Cells.Item(132, 11).Value = 4533
$ws.Cells.Item(132, 12).Value = 10098.9999
$ws.Cells.Item(132, 13).Value = -2003
$ws.Cells.Item(132, 14).Value = -15158.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2998.0908
$ws.Cells.Item(81, 9).Value = 3000
$ws.Cells.Item(81, 10).Value = 2997.9
$ws.Cells.Item(81, 11).Value = 6000
$ws.Cells.Item(81, 12).Value = 5995.8
$ws.Cells.Item(81, 13).Value = -4939
$ws.Cells.Item(81, 14).Value = -8117.8
$ws.Cells.Item(84, 8).Value = 2998.0908
$ws.Cells.Item(84, 9).Value = 3000
$ws.Cells.Item(84, 10).Value = 2997.9
$ws.Cells.Item(84, 11).Value = 30000
$ws.Cells.Item(84, 12).Value = 29979
$ws.Cells.Item(84, 13).Value = -24696
$ws.Cells.Item(84, 14).Value = -40587
$ws.Cells.Item(136, 8).Value = 1189.625
$ws.Cells.Item(136, 9).Value = 1059.6383
$ws.Cells.Item(136, 10).Value = 1868.4445
$ws.Cells.Item(136, 11).Value = 3178.9149
$ws.Cells.Item(136, 12).Value = 5605.333500000001
$ws.Cells.Item(136, 13).Value = -628.9149000000002
$ws.Cells.Item(136, 14).Value = -10705.3335
$ws.Cells.Item(140, 8).Value = 45582.637
$ws.Cells.Item(140, 10).Value = 45582.637
$ws.Cells.Item(140, 12).Value = 45582.637
$ws.Cells.Item(140, 14).Value = -55942.637
